$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that receive a value in each new row (row 18 has no D column, matching the target)
$rowCols = @{
    13 = @("A", "B", "C", "D", "E", "G", "I")
    14 = @("A", "B", "C", "D", "E", "G", "I")
    15 = @("A", "B", "C", "D", "E", "G", "I")
    16 = @("A", "B", "C", "D", "E", "G", "I")
    17 = @("A", "B", "C", "D", "E", "G", "I")
    18 = @("A", "B", "C", "E", "G", "I")
}

# Copy per-column formatting from the row above (row 12) for each new row 13-18
# so that style indices s="4" (A/B) and s="2" (C-I) are reused without new styles being created.
# Only the columns that will actually receive a value in each row are touched, so no
# stray empty styled cells are introduced (e.g. row 18 has no D column).
for ($r = 13; $r -le 18; $r++) {
    $targetCols = $rowCols[$r]
    foreach ($col in $targetCols) {
        $ws.Range(($col + "12")).Copy()
        $ws.Range(($col + $r)).PasteSpecial(-4122) | Out-Null
    }
    # Match the workbook default row height (28.3) explicitly, mirroring how the
    # pre-existing rows are normalized (ht="28.3" customHeight="1") by the engine.
    $ws.Rows($r).RowHeight = 28.3
}
$excel.CutCopyMode = 0

# Row 13
$ws.Range("A13").Value = "Disco.Localization.Resources"
$ws.Range("B13").Value = "Strings"
$ws.Range("C13").Value = "DateTime_just_now"
$ws.Range("D13").Value = "一分钟以内"
$ws.Range("E13").Value = "Just now"
$ws.Range("G13").Value = "Just now"
$ws.Range("I13").Value = "刚刚"

# Row 14
$ws.Range("A14").Value = "Disco.Localization.Resources"
$ws.Range("B14").Value = "Strings"
$ws.Range("C14").Value = "DateTime_several_minutes_ago"
$ws.Range("D14").Value = "60 分钟以内"
$ws.Range("E14").Value = "{0} minutes ago"
$ws.Range("G14").Value = "{0} minutes ago"
$ws.Range("I14").Value = "{0} 分钟前"

# Row 15
$ws.Range("A15").Value = "Disco.Localization.Resources"
$ws.Range("B15").Value = "Strings"
$ws.Range("C15").Value = "DateTime_several_hours_ago"
$ws.Range("D15").Value = "24 小时以内"
$ws.Range("E15").Value = "{0} hours ago"
$ws.Range("G15").Value = "{0} hours ago"
$ws.Range("I15").Value = "{0} 小时前"

# Row 16
$ws.Range("A16").Value = "Disco.Localization.Resources"
$ws.Range("B16").Value = "Strings"
$ws.Range("C16").Value = "DateTime_several_days_ago"
$ws.Range("D16").Value = "30 天以内"
$ws.Range("E16").Value = "{0} days ago"
$ws.Range("G16").Value = "{0} days ago"
$ws.Range("I16").Value = "{0} 天前"

# Row 17
$ws.Range("A17").Value = "Disco.Localization.Resources"
$ws.Range("B17").Value = "Strings"
$ws.Range("C17").Value = "DateTime_same_year"
$ws.Range("D17").Value = "同一年"
$ws.Range("E17").Value = "dd/MM"
$ws.Range("G17").Value = "dd/MM"
$ws.Range("I17").Value = "MM月dd日"

# Row 18
$ws.Range("A18").Value = "Disco.Localization.Resources"
$ws.Range("B18").Value = "Strings"
$ws.Range("C18").Value = "DateTime_date_only"
$ws.Range("E18").Value = "dd/MM/yyyy"
$ws.Range("G18").Value = "dd/MM/yyyy"
$ws.Range("I18").Value = "yyyy年MM月dd日"

# Update dimension implicitly handled by engine; set the active selection to C13 to match target
$ws.Range("C13").Select() | Out-Null

Write-Host "Edit applied successfully"
